# Applies the "456a3b4" gh-pages data refresh to 杭州-漫展信息.xlsx
#  - Sheet "展览"  (Exhibitions): refreshed "想去人数" (F) counters + one
#    refreshed cover image URL (I12).
#  - Sheet "演出"  (Performances): refreshed F11 counter, a brand-new show
#    inserted as row 17 ("《卡农》永恒经典名曲音乐会"), which pushes the
#    previously-existing rows 17-21 down to 18-22.
#  - Sheet "全部类型" (All types, union of every category): same refreshed
#    "想去人数" (F) counters + the same refreshed cover image URL (I12) as
#    "展览" (no row insertion needed there).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览": bump the "想去人数" (F) counters and refresh one cover URL.
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$expoCounts = @{
    3  = 120
    7  = 1532
    8  = 10949
    16 = 217
    18 = 236
    19 = 1166
    21 = 227
    22 = 712
    24 = 239
    26 = 701
    27 = 3320
    28 = 1033
    29 = 758
    33 = 950
    35 = 45
    37 = 6
    38 = 14
    39 = 1349
    40 = 3714
    41 = 5418
    43 = 104
    44 = 140
    45 = 238
    48 = 4084
}

foreach ($row in $expoCounts.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $expoCounts[$row]
}

$wsExpo.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202409/N8wn3kgd1727060970805.jpeg"

# ---------------------------------------------------------------------
# Sheet "演出": bump F11, insert the new row 17, and shift the rest down.
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Range("F11").Value = 479

# Insert a fresh row at position 17 - everything that used to live at
# rows 17-21 automatically slides down to rows 18-22.
$wsShow.Rows.Item(17).Insert()

# Give the new row 17 the same "序号" cell styling (border/bold/center)
# as every other row in column A, then fill in its sequence number.
$wsShow.Range("A16").Copy()
$wsShow.Range("A17").PasteSpecial(-4122)
$wsShow.Range("A17").Value = 16

# Force the date-looking value to stay plain text (matches every other
# "开始时间" cell in this column, which are plain strings, not Excel dates),
# then drop back to the default/no style so the cell's <c> tag stays bare.
$wsShow.Range("B17").NumberFormat = "@"
$wsShow.Range("B17").Value = "2024-12-20"
$wsShow.Range("B17").Style = "Normal"
$wsShow.Range("C17").Value = "杭州·【早鸟5折起】《卡农》永恒经典名曲音乐会"
$wsShow.Range("D17").Value = "曙光路31号 浙江音乐厅"
$wsShow.Range("E17").Value = "2024.12.20 19:30-12.20 21:00"
$wsShow.Range("F17").Value = 0
$wsShow.Range("G17").Value = 100
$wsShow.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=92724"
$wsShow.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202409/TzutCfwb1727056405062.jpeg"

# ---------------------------------------------------------------------
# Sheet "全部类型": same refreshed "想去人数" (F) counters + cover URL as
# "展览" - rows differ because this sheet merges every category together.
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$allCounts = @{
    9  = 10949
    15 = 217
    17 = 236
    18 = 1166
    20 = 227
    23 = 712
    25 = 239
    26 = 701
    27 = 1033
    30 = 758
    33 = 45
    34 = 6
    37 = 104
    38 = 140
    39 = 238
    43 = 4084
}

foreach ($row in $allCounts.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allCounts[$row]
}

$wsAll.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202409/N8wn3kgd1727060970805.jpeg"
